$d = $word.ActiveDocument

# 1. "From author annotations" -> "For author annotations"
$d.Content.Find.Execute("From author annotations", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "For author annotations", 2) | Out-Null

# 2. Insert two new bullet paragraphs ("Credit Title" and "Credit Authors") right after
#    the "Resource Title" paragraph, matching its numbering (ilvl=1, numId=1004).
$resourceTitle = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd("`r", "`a") -eq "Resource Title") {
        $resourceTitle = $p
        break
    }
}

$resourceTitle.Range.InsertParagraphAfter()
$creditTitlePara = $resourceTitle.Next()
$creditTitlePara.Range.Text = "Credit Title"

$creditTitlePara.Range.InsertParagraphAfter()
$creditAuthorsPara = $creditTitlePara.Next()
$creditAuthorsPara.Range.Text = "Credit Authors"
